$wb = $excel.ActiveWorkbook

# --- 1. "FCY current" sheet: just move the selection (no data changes) ---
$ws3 = $wb.Worksheets.Item("FCY current")
$ws3.Range("G1").Select() | Out-Null

# --- 2. "FCY Saving" -> rename to "LCSavingAccountInputter" and replace its
#        contents with the new 2x2 "Purpose of Bank Account / ID" table ---
$ws4 = $wb.Worksheets.Item("FCY Saving")
# Delete (not just clear) the old used range so formatting/row-height carried
# over from the previous data (e.g. the 16.5pt row) is removed as well.
$ws4.Range("A1:T2").Delete() | Out-Null

$ws4.Range("A1").Value = "Purpose of Bank Account"
$ws4.Range("A2").Value = "for property purpose use"
$ws4.Range("B1").Value = "ID"
$ws4.Range("B2").Value = 1007758835

$ws4.Name = "LCSavingAccountInputter"
$ws4.Columns.Item(2).ColumnWidth = 10.14
$ws4.Range("A2").Select() | Out-Null

# --- 3. New sheet "LCCurrentInputter", added after the last sheet.
#        Copy an existing plain sheet (no explicit pageSetup element) so the
#        new sheet keeps the same namespaces/formatting template, then wipe
#        its contents and fill in the new 2x2 "ID / Purpose of Bank Account"
#        table (mirrored layout vs. the sheet above). ---
$wsTemplate = $wb.Worksheets.Item("LCY Saving")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTemplate.Copy([System.Reflection.Missing]::Value, $lastSheet) | Out-Null

$ws5 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5.Name = "LCCurrentInputter"
$ws5.Range("A1:S2").Delete() | Out-Null

$ws5.Range("A1").Value = "ID"
$ws5.Range("B1").Value = "Purpose of Bank Account"
$ws5.Range("A2").Value = 1007414740
$ws5.Range("B2").Value = "for property purpose use"

$ws5.Columns.Item(1).ColumnWidth = 10.14
$ws5.Range("O8").Select() | Out-Null

Write-Host "Workbook updated: FCY Saving -> LCSavingAccountInputter, added LCCurrentInputter"
